$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '  -2.46%  '
$ws.Range("E3").Value = '  -5.11%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  -1.53%  '
$ws.Range("E6").Value = '  -7.40%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -11.75%  '
$ws.Range("E9").Value = '  -5.17%  '
$ws.Range("E10").Value = '  -2.86%  '
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("E12").Value = '  -3.29%  '
$ws.Range("E13").Value = '  -3.77%  '
$ws.Range("E14").Value = '  -6.83%  '
$ws.Range("D15").Value = '2.761.73'
$ws.Range("E15").Value = '  -4.51%  '
$ws.Range("E16").Value = '  -2.64%  '
$ws.Range("E17").Value = '  -6.67%  '
$ws.Range("D18").Value = '2.346.77'
$ws.Range("E18").Value = '  -4.18%  '
$ws.Range("E19").Value = '  -3.24%  '
$ws.Range("E20").Value = '  -1.94%  '
$ws.Range("E21").Value = '  -2.11%  '
$ws.Range("E22").Value = '  -6.99%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  -2.10%  '
$ws.Range("E25").Value = '  -11.61%  '
$ws.Range("E26").Value = '  +4.57%  '
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("E28").Value = '  -5.05%  '
$ws.Range("E29").Value = '  -5.24%  '
$ws.Range("E30").Value = '  -8.99%  '
$ws.Range("E31").Value = '  -13.02%  '
$ws.Range("E32").Value = '  -2.46%  '
$ws.Range("E33").Value = '  -9.67%  '
$ws.Range("E34").Value = '  -5.79%  '
$ws.Range("E35").Value = '  -6.82%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("E40").Value = '  -10.08%  '
$ws.Range("E41").Value = '  +1.45%  '
$ws.Range("E42").Value = '  -1.37%  '
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("E44").Value = '  +0.06%  '
$c = $ws.Range("D45")
$c.Value2 = "'140.15"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.76%  '
$ws.Range("E46").Value = '  -2.85%  '
$ws.Range("E47").Value = '  -10.87%  '
$ws.Range("E48").Value = '  -4.14%  '
$ws.Range("E49").Value = '  -4.16%  '
$c = $ws.Range("D50")
$c.Value2 = "'18.81"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -11.63%  '
$ws.Range("E51").Value = '  -4.85%  '
